$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Flip the "NO" markers to "YES" for rows that now have topper/percentile data (24-28)
$ws.Range("S24").Value = "YES"
$ws.Range("S25").Value = "YES"
$ws.Range("S26").Value = "YES"
$ws.Range("S27").Value = "YES"
$ws.Range("S28").Value = "YES"

# 2. Fill in the Topper column (N) for the two "GO CLASSES MOCK" rows that
#    didn't have a topper score recorded yet.
$ws.Range("N31").Value = 81
$ws.Range("N32").Value = 81

# 3. Add the new "GO CLASSES MOCK - 4" test row (row 33)
$ws.Range("A33").Value = 44925
$ws.Range("B33").Value = 32
$ws.Range("C33").Value = "GO CLASSES MOCK - 4"
$ws.Range("E33").Value = 58
$ws.Range("F33").Value = 48
$ws.Range("G33").Value = 10
$ws.Range("H33").Value = 7
$ws.Range("I33").Value = 71
$ws.Range("J33").Value = -1
$ws.Range("K33").Value = 28
$ws.Range("L33").Value = 70
$ws.Range("M33").Value = 44.97
$ws.Range("N33").Value = 79
$ws.Range("P33").Value = 6
$ws.Range("Q33").Value = 68
$ws.Range("S31").Value = "YES"
$ws.Range("S32").Value = "YES"
$ws.Range("S33").Value = "YES"

# R, T columns use fill-down formulas relative to the preceding rows
$ws.Range("R33").Formula = "=(Q33-P33+1)/(Q33)*100"
$ws.Range("T31").Formula = "=N31-L31"
$ws.Range("T32").Formula = "=N32-L32"
$ws.Range("T33").Formula = "=N33-L33"
$ws.Range("T34").Formula = "=N34-L34"

# 4. Restore the view state recorded after the edit (scrolled down, new selection)
$excel.ActiveWindow.ScrollRow = 44
$ws.Range("L30").Select()
